# "Generate Report for Handoff" - refresh the localization-status report:
# the zh-cn / de-de rows move from "In Translation" to "Ready for handoff"
# and pick up the new handoff generation timestamps, and the
# Status / Generate-Date columns are widened to fit the longer text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E, F) and the
# "Latest HO Xliff Generate Date" column (G)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-18 18:41:26"

# zh-cn detail sheet: Status (C) and Latest Handoff Datetime (H)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-18 18:41:20"

# de-de detail sheet: Status (C) and Latest Handoff Datetime (H)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-18 18:41:26"

# Widen the Status / Generate-Date columns so the new text fits
# (engine rounds ColumnWidth to its internal pixel grid, so back the
# character width off slightly to land as close as possible to the
# target column width after that rounding)
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333

$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
